$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: fill Image (E), Description (D), Name (A) and Price (C) per row
# (matches the shared-string insertion order seen in the target workbook)
$ws.Range("E17").Style = "Normal"
$ws.Range("E17").Value = "https://i.postimg.cc/dtG6DkJj/Whats-App-Image-2025-05-28-at-14-10-40.jpg"
$ws.Range("D17").Value = "UNDER 19 - 32-38 - 3PCS BOX - NS LYCRA - 235"
$ws.Range("A17").Value = "UNDER 19"
$ws.Range("C17").Style = "Normal"
$ws.Range("C17").Value = 235
$ws.Range("E18").Style = "Normal"
$ws.Range("E18").Value = "https://i.postimg.cc/KckJSpVV/Whats-App-Image-2025-05-28-at-14-10-40-1.jpg"
$ws.Range("D18").Value = "RICH BOND/A-3281 - 24-38 - 3PCS BOX - WASHING - 565"
$ws.Range("A18").Value = "RICH BOND/A-3281"
$ws.Range("C18").Style = "Normal"
$ws.Range("C18").Value = 565
$ws.Range("E19").Style = "Normal"
$ws.Range("E19").Value = "https://i.postimg.cc/wMmVX3LH/Whats-App-Image-2025-05-28-at-14-10-41.jpg"
$ws.Range("D19").Value = "RICH BOND/A-219 - 24-38 - 3PCS BOX - JORDEN CRUSH - 310"
$ws.Range("A19").Value = "RICH BOND/A-219"
$ws.Range("C19").Style = "Normal"
$ws.Range("C19").Value = 310
$ws.Range("E20").Style = "Normal"
$ws.Range("E20").Value = "https://i.postimg.cc/YCk3Fr9K/Whats-App-Image-2025-05-28-at-14-10-41-1.jpg"
$ws.Range("D20").Value = "RELAX WEAR/A-173 - 30-36 - 3PCS BOX - PC/COTTON - 220"
$ws.Range("A20").Value = "RELAX WEAR/A-173"
$ws.Range("C20").Style = "Normal"
$ws.Range("C20").Value = 220
$ws.Range("E21").Style = "Normal"
$ws.Range("E21").Value = "https://i.postimg.cc/xCkyQqsH/Whats-App-Image-2025-05-28-at-14-10-42.jpg"
$ws.Range("D21").Value = "RETRO/A-715 - 24-38 - 3PCS BOX - COTTON - 225"
$ws.Range("A21").Value = "RETRO/A-715"
$ws.Range("C21").Style = "Normal"
$ws.Range("C21").Value = 225
$ws.Range("E22").Style = "Normal"
$ws.Range("E22").Value = "https://i.postimg.cc/GpgJ6vjd/Whats-App-Image-2025-05-28-at-14-10-42-1.jpg"
$ws.Range("D22").Value = "RETRO/A-712 - 24-28 - 3PCS BOX - COTTON - 225"
$ws.Range("A22").Value = "RETRO/A-712"
$ws.Range("C22").Style = "Normal"
$ws.Range("C22").Value = 225
$ws.Range("E23").Style = "Normal"
$ws.Range("E23").Value = "https://i.postimg.cc/Nf88B3xN/Whats-App-Image-2025-05-28-at-14-10-43.jpg"
$ws.Range("D23").Value = "CATTNOS/A-2002 - 24-36 - 6PCS BOX - PC/COTTON - 200"
$ws.Range("A23").Value = "CATTNOS/A-2002"
$ws.Range("C23").Style = "Normal"
$ws.Range("C23").Value = 200
$ws.Range("E24").Style = "Normal"
$ws.Range("E24").Value = "https://i.postimg.cc/xCWG1Y7Z/Whats-App-Image-2025-05-28-at-14-10-43-1.jpg"
$ws.Range("D24").Value = "CATTNOS/A-2001 - 30-36 - 6PCS BOX - PC/COTTON - 215"
$ws.Range("A24").Value = "CATTNOS/A-2001"
$ws.Range("C24").Style = "Normal"
$ws.Range("C24").Value = 215
$ws.Range("E25").Style = "Normal"
$ws.Range("E25").Value = "https://i.postimg.cc/ryrCY0R2/Whats-App-Image-2025-05-28-at-14-10-43-2.jpg"
$ws.Range("D25").Value = "CATTNOS/A-2007 - 24-36 - 6PCS BOX - PC/COTTON - 195"
$ws.Range("A25").Value = "CATTNOS/A-2007"
$ws.Range("C25").Style = "Normal"
$ws.Range("C25").Value = 195
$ws.Range("E26").Style = "Normal"
$ws.Range("E26").Value = "https://i.postimg.cc/KYCx9sCb/Whats-App-Image-2025-05-28-at-14-10-44.jpg"
$ws.Range("D26").Value = "RELAX WEAR/A-3301 - 26-36 - 3PCS BOX - DENIM - 330"
$ws.Range("A26").Value = "RELAX WEAR/A-3301"
$ws.Range("C26").Style = "Normal"
$ws.Range("C26").Value = 330
$ws.Range("E27").Style = "Normal"
$ws.Range("E27").Value = "https://i.postimg.cc/Wb7p28qm/Whats-App-Image-2025-05-28-at-14-10-44-1.jpg"
$ws.Range("D27").Value = "RETRO/A-203 - 24-38 - 3PCS BOX - COTTON - 180"
$ws.Range("A27").Value = "RETRO/A-203"
$ws.Range("C27").Style = "Normal"
$ws.Range("C27").Value = 180
$ws.Range("E28").Style = "Normal"
$ws.Range("E28").Value = "https://i.postimg.cc/NftfqSFf/Whats-App-Image-2025-05-28-at-14-10-45.jpg"
$ws.Range("D28").Value = "CORDIAL/A-6431 - 5----7 - 3PCS BOX - DENIM - 325"
$ws.Range("A28").Value = "CORDIAL/A-6431"
$ws.Range("C28").Style = "Normal"
$ws.Range("C28").Value = 325
$ws.Range("E29").Style = "Normal"
$ws.Range("E29").Value = "https://i.postimg.cc/WpgbyT2H/Whats-App-Image-2025-05-28-at-14-10-45-1.jpg"
$ws.Range("D29").Value = "CORDIAL/A-6434 - 5----7 - 3PCS BOX - DENIM - 350"
$ws.Range("A29").Value = "CORDIAL/A-6434"
$ws.Range("C29").Style = "Normal"
$ws.Range("C29").Value = 350
$ws.Range("E30").Style = "Normal"
$ws.Range("E30").Value = "https://i.postimg.cc/2SczW5wV/Whats-App-Image-2025-05-28-at-14-10-46.jpg"
$ws.Range("D30").Value = "CORDIAL/A-6434 - 11----13 - 3PCS BOX - DENIM - 410"
$ws.Range("A30").Value = "CORDIAL/A-6434"
$ws.Range("C30").Style = "Normal"
$ws.Range("C30").Value = 410
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").Value = "https://i.postimg.cc/t4Cq76FK/Whats-App-Image-2025-05-28-at-14-10-46-1.jpg"
$ws.Range("D31").Value = "CATTNOS/A-1039 - 24-36 - 3PCS BOX - COTTON - 150"
$ws.Range("A31").Value = "CATTNOS/A-1039"
$ws.Range("C31").Style = "Normal"
$ws.Range("C31").Value = 150
$ws.Range("E32").Style = "Normal"
$ws.Range("E32").Value = "https://i.postimg.cc/SRhq4Yyt/Whats-App-Image-2025-05-28-at-14-10-47.jpg"
$ws.Range("D32").Value = "CATTNOS/A-1033 - 24-36 - 3PCS BOX - COTTON - 145"
$ws.Range("A32").Value = "CATTNOS/A-1033"
$ws.Range("C32").Style = "Normal"
$ws.Range("C32").Value = 145
$ws.Range("E33").Style = "Normal"
$ws.Range("E33").Value = "https://i.postimg.cc/02H54f6V/Whats-App-Image-2025-05-28-at-14-10-55.jpg"
$ws.Range("D33").Value = "STELLER/A-2001 - 26-36 - 6PCS BOX - INTERLOCK - 270"
$ws.Range("A33").Value = "STELLER/A-2001"
$ws.Range("C33").Style = "Normal"
$ws.Range("C33").Value = 270
$ws.Range("E34").Style = "Normal"
$ws.Range("E34").Value = "https://i.postimg.cc/VkrsPCnm/Whats-App-Image-2025-05-28-at-14-10-56.jpg"
$ws.Range("D34").Value = "STELLER/A-2101 - 26-36 - 6PCS BOX - COTTON (JOGGERS) - 270"
$ws.Range("A34").Value = "STELLER/A-2101"
$ws.Range("C34").Style = "Normal"
$ws.Range("C34").Value = 270
$ws.Range("E35").Style = "Normal"
$ws.Range("E35").Value = "https://i.postimg.cc/xdL0cV6Z/Whats-App-Image-2025-05-28-at-14-10-56-1.jpg"
$ws.Range("D35").Value = "STELLER/A-2025 - 26-36 - 6PCS BOX - INTERLOCK - 280"
$ws.Range("A35").Value = "STELLER/A-2025"
$ws.Range("C35").Style = "Normal"
$ws.Range("C35").Value = 280
$ws.Range("E36").Style = "Normal"
$ws.Range("E36").Value = "https://i.postimg.cc/tTy96vf4/Whats-App-Image-2025-05-28-at-14-10-57.jpg"
$ws.Range("D36").Value = "STELLER/A-2002 - 26-36 - 6PCS BOX - INTERLOCK - 250"
$ws.Range("A36").Value = "STELLER/A-2002"
$ws.Range("C36").Style = "Normal"
$ws.Range("C36").Value = 250
$ws.Range("E37").Style = "Normal"
$ws.Range("E37").Value = "https://i.postimg.cc/qvPJJvxn/Whats-App-Image-2025-05-28-at-14-10-57-1.jpg"
$ws.Range("D37").Value = "STELLER/A-2103 - 26-36 - 6PCS BOX - COTTON ( STRAIGHT FIT) - 255"
$ws.Range("A37").Value = "STELLER/A-2103"
$ws.Range("C37").Style = "Normal"
$ws.Range("C37").Value = 255
$ws.Range("E38").Style = "Normal"
$ws.Range("E38").Value = "https://i.postimg.cc/Qx0tgDkK/Whats-App-Image-2025-05-28-at-14-56-30.jpg"
$ws.Range("D38").Value = "FABU LOVE/A-41110 - 20-32 - 3PCS BOX - COTTON LYCRA - 245"
$ws.Range("A38").Value = "FABU LOVE/A-41110"
$ws.Range("C38").Style = "Normal"
$ws.Range("C38").Value = 245
$ws.Range("E39").Style = "Normal"
$ws.Range("E39").Value = "https://i.postimg.cc/pLrdqSqG/Whats-App-Image-2025-05-28-at-14-56-30-1.jpg"
$ws.Range("D39").Value = "FABU LOVE/A-41239 - 20-32 - 3PCS BOX - COTTON LYCRA - 295"
$ws.Range("A39").Value = "FABU LOVE/A-41239"
$ws.Range("C39").Style = "Normal"
$ws.Range("C39").Value = 295
$ws.Range("E40").Style = "Normal"
$ws.Range("E40").Value = "https://i.postimg.cc/CLF5f5DV/Whats-App-Image-2025-05-28-at-14-56-31.jpg"
$ws.Range("D40").Value = "FABU LOVE/A-41487 - 22-32 - 3PCS BOX - COTTON LYCRA - 325"
$ws.Range("A40").Value = "FABU LOVE/A-41487"
$ws.Range("C40").Style = "Normal"
$ws.Range("C40").Value = 325
$ws.Range("E41").Style = "Normal"
$ws.Range("E41").Value = "https://i.postimg.cc/RC1qFKWR/Whats-App-Image-2025-05-28-at-14-56-31-1.jpg"
$ws.Range("D41").Value = "FABU LOVE/A-41295 - 22-38 - 3PCS BOX - COTTON LYCRA - 421"
$ws.Range("A41").Value = "FABU LOVE/A-41295"
$ws.Range("C41").Style = "Normal"
$ws.Range("C41").Value = 421
$ws.Range("E42").Style = "Normal"
$ws.Range("E42").Value = "https://i.postimg.cc/FH6Hqx69/Whats-App-Image-2025-05-28-at-14-56-32.jpg"
$ws.Range("D42").Value = "FABU LOVE/A-41214 - 22-36 - 3PCS BOX - COTTON LYCRA - 390"
$ws.Range("A42").Value = "FABU LOVE/A-41214"
$ws.Range("C42").Style = "Normal"
$ws.Range("C42").Value = 390
$ws.Range("E43").Style = "Normal"
$ws.Range("E43").Value = "https://i.postimg.cc/WbHDLjFL/Whats-App-Image-2025-05-28-at-14-56-33.jpg"
$ws.Range("D43").Value = "PINK FOX/A-1155 - 22-32 - 3PCS BOX - COTTON - 278"
$ws.Range("A43").Value = "PINK FOX/A-1155"
$ws.Range("C43").Style = "Normal"
$ws.Range("C43").Value = 278
$ws.Range("E44").Style = "Normal"
$ws.Range("E44").Value = "https://i.postimg.cc/RCG6dwv0/Whats-App-Image-2025-05-28-at-14-56-33-1.jpg"
$ws.Range("D44").Value = "PINK FOX/A-1957 - 22-32 - 3PCS BOX - COTTON - 238"
$ws.Range("A44").Value = "PINK FOX/A-1957"
$ws.Range("C44").Style = "Normal"
$ws.Range("C44").Value = 238
$ws.Range("E45").Style = "Normal"
$ws.Range("E45").Value = "https://i.postimg.cc/vHB1NDmt/Whats-App-Image-2025-05-28-at-14-56-33-2.jpg"
$ws.Range("D45").Value = "PINK FOX/A-1918 - 22-36 - 3PCS BOX - COTTON - 253"
$ws.Range("A45").Value = "PINK FOX/A-1918"
$ws.Range("C45").Style = "Normal"
$ws.Range("C45").Value = 253
$ws.Range("E46").Style = "Normal"
$ws.Range("E46").Value = "https://i.postimg.cc/3wMdNBFt/Whats-App-Image-2025-05-28-at-14-56-34.jpg"
$ws.Range("D46").Value = "FABU LOVE/A-41275 - 22-38 - 3PCS BOX - COTTON LYCRA - 340"
$ws.Range("A46").Value = "FABU LOVE/A-41275"
$ws.Range("C46").Style = "Normal"
$ws.Range("C46").Value = 340
$ws.Range("E47").Style = "Normal"
$ws.Range("E47").Value = "https://i.postimg.cc/wjR19PLt/Whats-App-Image-2025-05-28-at-14-56-36.jpg"
$ws.Range("D47").Value = "PINK FOX/A-1784 - 22-36 - 3PCS BOX - TENCIL - 335"
$ws.Range("A47").Value = "PINK FOX/A-1784"
$ws.Range("C47").Style = "Normal"
$ws.Range("C47").Value = 335

# Phase 2: fill Category (B) per row
$ws.Range("B17").Value = "LOWER"
$ws.Range("B18").Value = "LOWER"
$ws.Range("B19").Value = "LOWER"
$ws.Range("B20").Value = "LOWER"
$ws.Range("B21").Value = "LOWER"
$ws.Range("B22").Value = "LOWER"
$ws.Range("B23").Value = "CAPRI"
$ws.Range("B24").Value = "CAPRI"
$ws.Range("B25").Value = "CAPRI"
$ws.Range("B26").Value = "LOWER"
$ws.Range("B27").Value = "NICKER"
$ws.Range("B28").Value = "NICKER"
$ws.Range("B29").Value = "NICKER"
$ws.Range("B30").Value = "NICKER"
$ws.Range("B31").Value = "NICKER"
$ws.Range("B32").Value = "NICKER"
$ws.Range("B33").Value = "LOWER"
$ws.Range("B34").Value = "LOWER"
$ws.Range("B35").Value = "LOWER"
$ws.Range("B36").Value = "LOWER"
$ws.Range("B37").Value = "LOWER"
$ws.Range("B38").Value = "CORDSET"
$ws.Range("B39").Value = "CORDSET"
$ws.Range("B40").Value = "CORDSET"
$ws.Range("B41").Value = "CORDSET"
$ws.Range("B42").Value = "CORDSET"
$ws.Range("B43").Value = "CORDSET"
$ws.Range("B44").Value = "CORDSET"
$ws.Range("B45").Value = "CORDSET"
$ws.Range("B46").Value = "CORDSET"
$ws.Range("B47").Value = "CORDSET"

# Match the final selection left by the edit
$ws.Range("B17:B47").Select()
